$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B6").Value = "Ich denke es ist eine gute Idee, den Roboter zu verwenden."
